$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 2500  # H21: was 0
$ws.Cells.Item(21, 9).Value = 2500  # I21: was 0
$ws.Cells.Item(21, 11).Value = 2500  # K21: was 0
$ws.Cells.Item(21, 13).Value = -2032  # M21: was None

$ws.Cells.Item(23, 8).Value = 2500  # H23: was 0
$ws.Cells.Item(23, 9).Value = 2500  # I23: was 0
$ws.Cells.Item(23, 11).Value = 2500  # K23: was 0
$ws.Cells.Item(23, 13).Value = -2266  # M23: was None

$ws.Cells.Item(29, 8).Value = 100  # H29: was 75
$ws.Cells.Item(29, 10).Value = 0  # J29: was 50
$ws.Cells.Item(29, 12).Value = 0  # L29: was 150
$ws.Cells.Item(29, 14).ClearContents()  # N29: was -712

$ws.Cells.Item(53, 8).Value = 1011.86365  # H53: was 971.6087
$ws.Cells.Item(53, 9).Value = 832  # I53: was 785.375
$ws.Cells.Item(53, 11).Value = 832  # K53: was 785.375
$ws.Cells.Item(53, 13).Value = -195  # M53: was -148.375

$ws.Cells.Item(76, 8).Value = 8221.444  # H76: was 8624.875
$ws.Cells.Item(76, 9).Value = 4997.6665  # I76: was 4999.5
$ws.Cells.Item(76, 11).Value = 4997.6665  # K76: was 4999.5
$ws.Cells.Item(76, 13).Value = -4682.6665  # M76: was -4684.5

$ws.Cells.Item(79, 8).Value = 8221.444  # H79: was 8624.875
$ws.Cells.Item(79, 9).Value = 4997.6665  # I79: was 4999.5
$ws.Cells.Item(79, 11).Value = 4997.6665  # K79: was 4999.5
$ws.Cells.Item(79, 13).Value = -3905.6665  # M79: was -3907.5

$ws.Cells.Item(121, 8).Value = 1577.8572  # H121: was 1696.762
$ws.Cells.Item(121, 10).Value = 1599.2307  # J121: was 1738.5264
$ws.Cells.Item(121, 12).Value = 4797.6921  # L121: was 5215.5792
$ws.Cells.Item(121, 14).Value = -8291.6921  # N121: was -8709.5792

$ws.Cells.Item(137, 8).Value = 1531.3  # H137: was 1365.1428
$ws.Cells.Item(137, 9).Value = 1219.3334  # I137: was 1111.5
$ws.Cells.Item(137, 11).Value = 3658.0002  # K137: was 3334.5
$ws.Cells.Item(137, 13).Value = -1108.0002  # M137: was -784.5

$ws.Cells.Item(138, 8).Value = 3856.94  # H138: was 3801.889
$ws.Cells.Item(138, 9).Value = 2379.1  # I138: was 2361.4375
$ws.Cells.Item(138, 10).Value = 4490.3  # J138: was 4489.8657
$ws.Cells.Item(138, 11).Value = 7137.299999999999  # K138: was 7084.3125
$ws.Cells.Item(138, 12).Value = 13470.9  # L138: was 13469.5971
$ws.Cells.Item(138, 13).Value = -1997.299999999999  # M138: was -1944.3125
$ws.Cells.Item(138, 14).Value = -23750.9  # N138: was -23749.5971

$ws.Cells.Item(140, 8).Value = 199993.5  # H140: was 199994
$ws.Cells.Item(140, 10).Value = 199993.5  # J140: was 199994
$ws.Cells.Item(140, 12).Value = 199993.5  # L140: was 199994
$ws.Cells.Item(140, 14).Value = -210353.5  # N140: was -210354

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20347482  # H32: was 21145422
$ws.Cells.Item(32, 9).Value = 20732914  # I32: was 21615164
$ws.Cells.Item(32, 11).Value = 20732914  # K32: was 21615164
$ws.Cells.Item(32, 13).Value = -20732627  # M32: was -21614877

$ws.Cells.Item(61, 8).Value = 3104.647  # H61: was 2902.8223
$ws.Cells.Item(61, 9).Value = 2775.3076  # I61: was 2513.4375
$ws.Cells.Item(61, 10).Value = 4175  # J61: was 3861.3076
$ws.Cells.Item(61, 11).Value = 2775.3076  # K61: was 2513.4375
$ws.Cells.Item(61, 12).Value = 4175  # L61: was 3861.3076
$ws.Cells.Item(61, 13).Value = -2563.3076  # M61: was -2301.4375
$ws.Cells.Item(61, 14).Value = -4599  # N61: was -4285.3076

$ws.Cells.Item(63, 8).Value = 4355.875  # H63: was 4526.2666
$ws.Cells.Item(63, 9).Value = 2212.25  # I63: was 2271.1428
$ws.Cells.Item(63, 11).Value = 2212.25  # K63: was 2271.1428
$ws.Cells.Item(63, 13).Value = -1526.25  # M63: was -1585.1428

$ws.Cells.Item(66, 8).Value = 4355.875  # H66: was 4526.2666
$ws.Cells.Item(66, 9).Value = 2212.25  # I66: was 2271.1428
$ws.Cells.Item(66, 11).Value = 11061.25  # K66: was 11355.714
$ws.Cells.Item(66, 13).Value = -7629.25  # M66: was -7923.714

$ws.Cells.Item(109, 8).Value = 21750  # H109: was 28375
$ws.Cells.Item(109, 10).Value = 21750  # J109: was 28375
$ws.Cells.Item(109, 12).Value = 21750  # L109: was 28375
$ws.Cells.Item(109, 14).Value = -24524  # N109: was -31149

$ws.Cells.Item(110, 8).Value = 1564.35  # H110: was 1532.619
$ws.Cells.Item(110, 9).Value = 892.9375  # I110: was 893.2353000000001
$ws.Cells.Item(110, 11).Value = 892.9375  # K110: was 893.2353000000001
$ws.Cells.Item(110, 13).Value = 1152.0625  # M110: was 1151.7647

$ws.Cells.Item(122, 8).Value = 2814.9167  # H122: was 2902.6086
$ws.Cells.Item(122, 9).Value = 2528.2  # I122: was 2619.2632
$ws.Cells.Item(122, 11).Value = 7584.599999999999  # K122: was 7857.7896
$ws.Cells.Item(122, 13).Value = -5134.599999999999  # M122: was -5407.7896

$ws.Cells.Item(136, 8).Value = 3104.647  # H136: was 2902.8223
$ws.Cells.Item(136, 9).Value = 2775.3076  # I136: was 2513.4375
$ws.Cells.Item(136, 10).Value = 4175  # J136: was 3861.3076
$ws.Cells.Item(136, 11).Value = 8325.9228  # K136: was 7540.3125
$ws.Cells.Item(136, 12).Value = 12525  # L136: was 11583.9228
$ws.Cells.Item(136, 13).Value = -5775.9228  # M136: was -4990.3125
$ws.Cells.Item(136, 14).Value = -17625  # N136: was -16683.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 83333.336  # H132: was 92500
$ws.Cells.Item(132, 10).Value = 83333.336  # J132: was 92500
$ws.Cells.Item(132, 12).Value = 83333.336  # L132: was 92500
$ws.Cells.Item(132, 14).Value = -93453.336  # N132: was -102620

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 101796.5  # H20: was 103598
$ws.Cells.Item(20, 10).Value = 101796.5  # J20: was 103598
$ws.Cells.Item(20, 12).Value = 101796.5  # L20: was 103598
$ws.Cells.Item(20, 14).Value = -102268.5  # N20: was -104070

$ws.Cells.Item(30, 8).Value = 101796.5  # H30: was 103598
$ws.Cells.Item(30, 10).Value = 101796.5  # J30: was 103598
$ws.Cells.Item(30, 12).Value = 101796.5  # L30: was 103598
$ws.Cells.Item(30, 14).Value = -101978.5  # N30: was -103780

$ws.Cells.Item(31, 8).Value = 3846.9824  # H31: was 3941.4363
$ws.Cells.Item(31, 9).Value = 1741.75  # I31: was 1796.4445
$ws.Cells.Item(31, 11).Value = 1741.75  # K31: was 1796.4445
$ws.Cells.Item(31, 13).Value = -1446.75  # M31: was -1501.4445

$ws.Cells.Item(34, 8).Value = 3846.9824  # H34: was 3941.4363
$ws.Cells.Item(34, 9).Value = 1741.75  # I34: was 1796.4445
$ws.Cells.Item(34, 11).Value = 1741.75  # K34: was 1796.4445
$ws.Cells.Item(34, 13).Value = -1539.75  # M34: was -1594.4445

$ws.Cells.Item(105, 8).Value = 2550.5881  # H105: was 2647.6875
$ws.Cells.Item(105, 9).Value = 2739.5833  # I105: was 2898
$ws.Cells.Item(105, 11).Value = 2739.5833  # K105: was 2898
$ws.Cells.Item(105, 13).Value = -992.5832999999998  # M105: was -1151

$ws.Cells.Item(107, 8).Value = 925.6  # H107: was 48500.57
$ws.Cells.Item(107, 9).Value = 552.53845  # I107: was 71941.64
$ws.Cells.Item(107, 11).Value = 552.53845  # K107: was 71941.64
$ws.Cells.Item(107, 13).Value = 1367.46155  # M107: was -70021.64

$ws.Cells.Item(123, 8).Value = 105045.5  # H123: was 106117.336
$ws.Cells.Item(123, 10).Value = 104606.336  # J123: was 105757.25
$ws.Cells.Item(123, 12).Value = 104606.336  # L123: was 105757.25
$ws.Cells.Item(123, 14).Value = -114406.336  # N123: was -115557.25

$ws.Cells.Item(128, 8).Value = 101796.5  # H128: was 103598
$ws.Cells.Item(128, 10).Value = 101796.5  # J128: was 103598
$ws.Cells.Item(128, 12).Value = 101796.5  # L128: was 103598
$ws.Cells.Item(128, 14).Value = -111756.5  # N128: was -113558

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 62.9375  # H12: was 68.5
$ws.Cells.Item(12, 10).Value = 81.416664  # J12: was 92.90000000000001
$ws.Cells.Item(12, 12).Value = 244.249992  # L12: was 278.7
$ws.Cells.Item(12, 14).Value = -590.249992  # N12: was -624.7

$ws.Cells.Item(38, 8).Value = 77.35714  # H38: was 503.08334
$ws.Cells.Item(38, 9).Value = 65.666664  # I38: was 90.833336
$ws.Cells.Item(38, 10).Value = 98.40000000000001  # J38: was 915.3333
$ws.Cells.Item(38, 11).Value = 196.999992  # K38: was 272.500008
$ws.Cells.Item(38, 12).Value = 295.2  # L38: was 2745.9999
$ws.Cells.Item(38, 13).Value = 150.000008  # M38: was 74.49999200000002
$ws.Cells.Item(38, 14).Value = -989.2  # N38: was -3439.9999

$ws.Cells.Item(88, 8).Value = 3766.6667  # H88: was 3813.3333
$ws.Cells.Item(88, 10).Value = 3766.6667  # J88: was 3813.3333
$ws.Cells.Item(88, 12).Value = 11300.0001  # L88: was 11439.9999
$ws.Cells.Item(88, 14).Value = -12156.0001  # N88: was -12295.9999

$ws.Cells.Item(91, 8).Value = 3766.6667  # H91: was 3813.3333
$ws.Cells.Item(91, 10).Value = 3766.6667  # J91: was 3813.3333
$ws.Cells.Item(91, 12).Value = 11300.0001  # L91: was 11439.9999
$ws.Cells.Item(91, 14).Value = -14264.0001  # N91: was -14403.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4500  # H70: was 4484.8184
$ws.Cells.Item(70, 9).Value = 0  # I70: was 4500
$ws.Cells.Item(70, 10).Value = 4500  # J70: was 4484.3438
$ws.Cells.Item(70, 11).Value = 0  # K70: was 4500
$ws.Cells.Item(70, 12).Value = 4500  # L70: was 4484.3438
$ws.Cells.Item(70, 13).ClearContents()  # M70: was -4230
$ws.Cells.Item(70, 14).Value = -5040  # N70: was -5024.3438

$ws.Cells.Item(73, 8).Value = 4500  # H73: was 4484.8184
$ws.Cells.Item(73, 9).Value = 0  # I73: was 4500
$ws.Cells.Item(73, 10).Value = 4500  # J73: was 4484.3438
$ws.Cells.Item(73, 11).Value = 0  # K73: was 4500
$ws.Cells.Item(73, 12).Value = 4500  # L73: was 4484.3438
$ws.Cells.Item(73, 13).ClearContents()  # M73: was -3564
$ws.Cells.Item(73, 14).Value = -6372  # N73: was -6356.3438

$ws.Cells.Item(82, 8).Value = 132000  # H82: was 130000
$ws.Cells.Item(82, 10).Value = 132000  # J82: was 130000
$ws.Cells.Item(82, 12).Value = 132000  # L82: was 130000
$ws.Cells.Item(82, 14).Value = -132766  # N82: was -130766

$ws.Cells.Item(85, 8).Value = 132000  # H85: was 130000
$ws.Cells.Item(85, 10).Value = 132000  # J85: was 130000
$ws.Cells.Item(85, 12).Value = 132000  # L85: was 130000
$ws.Cells.Item(85, 14).Value = -134652  # N85: was -132652

$ws.Cells.Item(102, 8).Value = 2777.2273  # H102: was 2818.524
$ws.Cells.Item(102, 9).Value = 2630.25  # I102: was 2668.158
$ws.Cells.Item(102, 11).Value = 2630.25  # K102: was 2668.158
$ws.Cells.Item(102, 13).Value = -1008.25  # M102: was -1046.158

$ws.Cells.Item(107, 8).Value = 1075.4667  # H107: was 1078.6923
$ws.Cells.Item(107, 10).Value = 1720.5714  # J107: was 1987
$ws.Cells.Item(107, 12).Value = 1720.5714  # L107: was 1987
$ws.Cells.Item(107, 14).Value = -5560.5714  # N107: was -5827

$ws.Cells.Item(113, 8).Value = 39541.39  # H113: was 46663.316
$ws.Cells.Item(113, 9).Value = 6941  # I113: was 8409.5
$ws.Cells.Item(113, 10).Value = 53804.062  # J113: was 56864.332
$ws.Cells.Item(113, 11).Value = 6941  # K113: was 8409.5
$ws.Cells.Item(113, 12).Value = 53804.062  # L113: was 56864.332
$ws.Cells.Item(113, 13).Value = -4771  # M113: was -6239.5
$ws.Cells.Item(113, 14).Value = -58144.062  # N113: was -61204.332

$ws.Cells.Item(125, 8).Value = 111117.6  # H125: was 111323.4
$ws.Cells.Item(125, 10).Value = 111117.6  # J125: was 111323.4
$ws.Cells.Item(125, 12).Value = 111117.6  # L125: was 111323.4
$ws.Cells.Item(125, 14).Value = -116037.6  # N125: was -116243.4

$ws.Cells.Item(132, 8).Value = 2970.8262  # H132: was 3022
$ws.Cells.Item(132, 9).Value = 2726.5  # I132: was 2772.8948
$ws.Cells.Item(132, 11).Value = 8179.5  # K132: was 8318.6844
$ws.Cells.Item(132, 13).Value = -5649.5  # M132: was -5788.6844

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9318.706  # H7: was 9894.267
$ws.Cells.Item(7, 9).Value = 7490.6665  # I7: was 8201.714
$ws.Cells.Item(7, 11).Value = 7490.6665  # K7: was 8201.714
$ws.Cells.Item(7, 13).Value = -7378.6665  # M7: was -8089.714

$ws.Cells.Item(46, 8).Value = 2455.5  # H46: was 2842.5
$ws.Cells.Item(46, 9).Value = 1658.4  # I46: was 1723
$ws.Cells.Item(46, 10).Value = 2628.7827  # J46: was 3122.375
$ws.Cells.Item(46, 11).Value = 1658.4  # K46: was 1723
$ws.Cells.Item(46, 12).Value = 2628.7827  # L46: was 3122.375
$ws.Cells.Item(46, 13).Value = -1470.4  # M46: was -1535
$ws.Cells.Item(46, 14).Value = -3004.7827  # N46: was -3498.375

$ws.Cells.Item(61, 8).Value = 18498.555  # H61: was 18498.777
$ws.Cells.Item(61, 9).Value = 16641  # I61: was 16641.285
$ws.Cells.Item(61, 11).Value = 16641  # K61: was 16641.285
$ws.Cells.Item(61, 13).Value = -16439  # M61: was -16439.285

$ws.Cells.Item(113, 8).Value = 18498.555  # H113: was 18498.777
$ws.Cells.Item(113, 9).Value = 16641  # I113: was 16641.285
$ws.Cells.Item(113, 11).Value = 16641  # K113: was 16641.285
$ws.Cells.Item(113, 13).Value = -14471  # M113: was -14471.285

$ws.Cells.Item(126, 8).Value = 9318.706  # H126: was 9894.267
$ws.Cells.Item(126, 9).Value = 7490.6665  # I126: was 8201.714
$ws.Cells.Item(126, 11).Value = 22471.9995  # K126: was 24605.142
$ws.Cells.Item(126, 13).Value = -20001.9995  # M126: was -22135.142

$ws.Cells.Item(141, 8).Value = 150000  # H141: was 383333.34
$ws.Cells.Item(141, 10).Value = 150000  # J141: was 383333.34
$ws.Cells.Item(141, 12).Value = 150000  # L141: was 383333.34
$ws.Cells.Item(141, 14).Value = -160360  # N141: was -393693.34

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 109994.5  # H119: was 109999
$ws.Cells.Item(119, 10).Value = 109994.5  # J119: was 109999
$ws.Cells.Item(119, 12).Value = 109994.5  # L119: was 109999
$ws.Cells.Item(119, 14).Value = -119670.5  # N119: was -119675

$ws.Cells.Item(126, 8).Value = 4331.8335  # H126: was 4332
$ws.Cells.Item(126, 9).Value = 4331.8335  # I126: was 4332
$ws.Cells.Item(126, 11).Value = 12995.5005  # K126: was 12996
$ws.Cells.Item(126, 13).Value = -10525.5005  # M126: was -10526

$ws.Cells.Item(128, 8).Value = 0  # H128: was 81220
$ws.Cells.Item(128, 10).Value = 0  # J128: was 81220
$ws.Cells.Item(128, 12).Value = 0  # L128: was 81220
$ws.Cells.Item(128, 14).ClearContents()  # N128: was -91180

$ws.Cells.Item(132, 8).Value = 29298.5  # H132: was 24068.426
$ws.Cells.Item(132, 9).Value = 31546.943  # I132: was 24978.133
$ws.Cells.Item(132, 10).Value = 3066.6667  # J132: was 3600
$ws.Cells.Item(132, 11).Value = 94640.829  # K132: was 74934.399
$ws.Cells.Item(132, 12).Value = 9200.000100000001  # L132: was 10800
$ws.Cells.Item(132, 13).Value = -92110.829  # M132: was -72404.399
$ws.Cells.Item(132, 14).Value = -14260.0001  # N132: was -15860
